# Wind renamed the steel "Wire Rod: Operating Rate of Main Steel Plant"
# indicator; the plot's source sheet now needs a dedicated row carrying
# the new Wind label ("Operating Rate: Wire Rod") paired with the
# existing Chinese label, inserted right after the old wire-rod row.
# Every row below shifts down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(115).EntireRow.Insert()

$ws.Range("A115").Value = "Operating Rate: Wire Rod"
$ws.Range("B115").Value = "线材：主要钢厂开工率"
